# Expense Report update:
#  - Refund note on row 22 (Amazon Order 9) updated to include the new $999 power-adapter refund
#  - Two new expense rows added (Amazon Order 13 / Amazon Order 14) for the Tower modification
#  - Bill-of-materials text for the previously-listed power adapter row clarified ("Test")
#  - Totals formula extended to include the new rows
#  - Footnote text at the bottom of the sheet updated with the new combined refund note

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the refund formula in row 22 to subtract the additional $999 refund
$ws.Range("E22").Formula = "=1537.58-8.18-326-999"

# 2. Insert two new rows right before the current "Total Expenses" row (old row 29)
$ws.Rows("29:30").Insert()

# 3. Re-label the existing row 28 items text: the "12V Power Adapters" BOM line is now the
#    "test" batch, since a "finalized" batch is being added for the new orders below.
$ws.Range("F28").Value = "12V Test Power Adapters"

# 4. Fill in new row 29 - Amazon Order 13 (Dr. Feron)
#    (Rows.Insert() already copied down the date/hyperlink/currency styles from row 28,
#    so simply assigning the values keeps the correct number formats. Hyperlinks.Add()
#    re-applies its own "Hyperlink" style variant, so re-apply the named style afterwards
#    to keep using the workbook's existing Hyperlink cell style.)
$ws.Range("A29").Value = "Amazon Order 13"
$ws.Range("B29").Value = 43164
$ws.Range("C29").Value = "Dr. Feron"
$ws.Range("D29").Value = "Amazon Order 13.pdf"
$ws.Range("E29").Value = 995.37
$ws.Range("F29").Value = "12V Finalized Power Adapters"
$ws.Hyperlinks.Add($ws.Range("D29"), "Amazon%20Order%2013.pdf")
$ws.Range("D29").Style = "Hyperlink"

# 5. Fill in new row 30 - Amazon Order 14 (Morgan)
$ws.Range("A30").Value = "Amazon Order 14"
$ws.Range("B30").Value = 43170
$ws.Range("C30").Value = "Morgan"
$ws.Range("D30").Value = "Amazon Order 14.pdf"
$ws.Range("E30").Value = 25.01
$ws.Range("F30").Value = "M4 screws & nuts"
$ws.Hyperlinks.Add($ws.Range("D30"), "Amazon%20Order%2014.pdf")
$ws.Range("D30").Style = "Hyperlink"

# 6. Extend the "Total Expenses" SUM formula (now on row 31) to cover the new rows
$ws.Range("E31").Formula = "=SUM(E2:E30)"

# 7. Update the bottom footnote (now row 33) with the combined refund text
$ws.Range("D33").Value = "***Refunded $8.18 for an incorrect shipment, an additional $326.00 for returned USB cables, and another $999.00 for returned power adapters"

# 8. Update selection to match the saved view state
$ws.Range("D35").Select()
